$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name and title text
$ws.Name = "Through 2022-10-24"
$ws.Range("A11").Value = "October (through 10-24)"

# Update October row (row 11)
$ws.Range("B11").Value = 21
$ws.Range("C11").Value = 40
$ws.Range("D11").Value = 54
$ws.Range("F11").Value = 41
$ws.Range("G11").Value = 120
$ws.Range("H11").Value = 151
$ws.Range("I11").Value = 89

# Update Total row (row 12)
$ws.Range("B12").Value = 247
$ws.Range("C12").Value = 469
$ws.Range("D12").Value = 681
$ws.Range("F12").Value = 463
$ws.Range("G12").Value = 1021
$ws.Range("H12").Value = 1398
$ws.Range("I12").Value = 1366
